# Correct the casing of three label cells in column A (Sheet1):
#   A2: mdaTextHomepage -> mdaTextHomePage
#   A4: MdaTitle        -> mdaTitle
#   A8: pageTitlenewTab -> pageTitleNewTab
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "mdaTextHomePage"
$ws.Range("A4").Value = "mdaTitle"
$ws.Range("A8").Value = "pageTitleNewTab"

# Move the active selection from B3 to A2, matching the saved view state.
$ws.Range("A2").Select()
